# Daily attendance processing - 2026-01-11 12:51:48
# Swap the "Recorded By" (column G) listing order from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every session row in the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$recordedByCol = 7   # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
